$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.4
$ws.Range("C2").Value = 2.15
$ws.Range("D2").Value = 0.96
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.9
$ws.Range("G2").Value = 0.24
$ws.Range("H2").Value = 2.23
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0.06

# Row 3
$ws.Range("B3").Value = 1.44
$ws.Range("D3").Value = 2.19
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.9
$ws.Range("G3").Value = 1.22
$ws.Range("I3").Value = 2.8
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0.08

# Row 4
$ws.Range("B4").Value = 1.04
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.9
$ws.Range("G4").Value = 1.24
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 0.99
$ws.Range("K4").Value = 0.1
